$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos data (prices / volume %) scraped on Mon Aug 28 2023

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.248.01'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.657.31'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.76%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.49'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5240'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.14%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.69%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2669'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06351'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.69'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.68%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.557'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.658.33'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.884.76'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5679'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8145'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.59'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.233.21'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.725'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.62'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.35'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.042'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.60'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1203'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.287'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.00'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.496'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05624'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.280'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.511'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.384'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.585'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9477'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.39%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5780'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.93%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.908'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8480'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.26%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.033.39'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.32'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '58.58'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₈107'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.001'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.22%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05320'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.85%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.040'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.27%  '
